$d = $word.ActiveDocument

$d.Content.Find.Execute("420×2=840", $true, $false, $false, $false, $false, $true, 1, $false, "996×3=2988", 2) | Out-Null
$d.Content.Find.Execute("663×9=5967", $true, $false, $false, $false, $false, $true, 1, $false, "852×4=3408", 2) | Out-Null
$d.Content.Find.Execute("189×2=378", $true, $false, $false, $false, $false, $true, 1, $false, "357×2=714", 2) | Out-Null
$d.Content.Find.Execute("985×7=6895", $true, $false, $false, $false, $false, $true, 1, $false, "848×8=6784", 2) | Out-Null
$d.Content.Find.Execute("711×5=3555", $true, $false, $false, $false, $false, $true, 1, $false, "280×9=2520", 2) | Out-Null
$d.Content.Find.Execute("924×3=2772", $true, $false, $false, $false, $false, $true, 1, $false, "954×5=4770", 2) | Out-Null
$d.Content.Find.Execute("707×6=4242", $true, $false, $false, $false, $false, $true, 1, $false, "312×8=2496", 2) | Out-Null
$d.Content.Find.Execute("893×9=8037", $true, $false, $false, $false, $false, $true, 1, $false, "943×2=1886", 2) | Out-Null
$d.Content.Find.Execute("483×5=2415", $true, $false, $false, $false, $false, $true, 1, $false, "401×5=2005", 2) | Out-Null
$d.Content.Find.Execute("550×3=1650", $true, $false, $false, $false, $false, $true, 1, $false, "967×8=7736", 2) | Out-Null
$d.Content.Find.Execute("248×5=1240", $true, $false, $false, $false, $false, $true, 1, $false, "992×6=5952", 2) | Out-Null
$d.Content.Find.Execute("652×9=5868", $true, $false, $false, $false, $false, $true, 1, $false, "855×2=1710", 2) | Out-Null
$d.Content.Find.Execute("690×7=4830", $true, $false, $false, $false, $false, $true, 1, $false, "599×4=2396", 2) | Out-Null
$d.Content.Find.Execute("891×4=3564", $true, $false, $false, $false, $false, $true, 1, $false, "924×4=3696", 2) | Out-Null
$d.Content.Find.Execute("660×9=5940", $true, $false, $false, $false, $false, $true, 1, $false, "451×8=3608", 2) | Out-Null
$d.Content.Find.Execute("915×2=1830", $true, $false, $false, $false, $false, $true, 1, $false, "558×5=2790", 2) | Out-Null
$d.Content.Find.Execute("698×8=5584", $true, $false, $false, $false, $false, $true, 1, $false, "712×4=2848", 2) | Out-Null
$d.Content.Find.Execute("836×4=3344", $true, $false, $false, $false, $false, $true, 1, $false, "258×3=774", 2) | Out-Null
$d.Content.Find.Execute("771×5=3855", $true, $false, $false, $false, $false, $true, 1, $false, "305×8=2440", 2) | Out-Null
$d.Content.Find.Execute("390×3=1170", $true, $false, $false, $false, $false, $true, 1, $false, "479×9=4311", 2) | Out-Null
$d.Content.Find.Execute("456×8=3648", $true, $false, $false, $false, $false, $true, 1, $false, "419×6=2514", 2) | Out-Null
$d.Content.Find.Execute("194×4=776", $true, $false, $false, $false, $false, $true, 1, $false, "871×2=1742", 2) | Out-Null
$d.Content.Find.Execute("435×4=1740", $true, $false, $false, $false, $false, $true, 1, $false, "208×8=1664", 2) | Out-Null
$d.Content.Find.Execute("320×3=960", $true, $false, $false, $false, $false, $true, 1, $false, "886×5=4430", 2) | Out-Null
$d.Content.Find.Execute("991×6=5946", $true, $false, $false, $false, $false, $true, 1, $false, "755×7=5285", 2) | Out-Null
